$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row right after the header (new row 2), shifting all
# existing data rows down by one. Clear any inherited formatting so the
# new row stays style-less like the other data rows.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()
$ws.Range("A2").Value = 1.779657483100891
$ws.Range("B2").Value = 1.38677453994751
$ws.Range("C2").Value = -1.217397570610046

# Append 9 new data rows after the (now shifted) last existing row (row 22).
$ws.Range("A23").Value = 0.4938832223415375
$ws.Range("B23").Value = -9.879995346069336
$ws.Range("C23").Value = -0.3310975134372711

$ws.Range("A24").Value = 6.862566947937012
$ws.Range("B24").Value = 17.44670104980469
$ws.Range("C24").Value = -6.492054462432861

$ws.Range("A25").Value = 0.584963321685791
$ws.Range("B25").Value = -0.4454802870750427
$ws.Range("C25").Value = 2.492385864257812

$ws.Range("A26").Value = -7.439141273498535
$ws.Range("B26").Value = -3.388273239135742
$ws.Range("C26").Value = 10.04511070251465

$ws.Range("A27").Value = 1.835583806037903
$ws.Range("B27").Value = 12.65514183044434
$ws.Range("C27").Value = -3.548196077346802

$ws.Range("A28").Value = -0.2848250865936279
$ws.Range("B28").Value = -6.074337482452393
$ws.Range("C28").Value = 3.493201971054077

$ws.Range("A29").Value = 7.885753154754639
$ws.Range("B29").Value = -7.940680980682373
$ws.Range("C29").Value = -2.813696384429932

$ws.Range("A30").Value = 1.459545493125916
$ws.Range("B30").Value = 0.23096264898777
$ws.Range("C30").Value = -2.083990097045898

$ws.Range("A31").Value = -3.29206657409668
$ws.Range("B31").Value = 4.709334373474121
$ws.Range("C31").Value = 3.162437200546265
